# CN. Maximum Path sum (recurrssion,tabulation,memoization,space optimization)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate row 9 with the new "Maximum Path Sum in the matrix" entry
$ws.Range("A9").Value = "CN"
$ws.Range("B9").Value = "CN"
$ws.Range("C9").Value = "Maximum Path Sum in the matrix"
$ws.Range("D9").Value = "Java"
$ws.Range("E9").Value = "DP(Recurrsion+Memonization+Tabulation+space optimization)"

# Match the formatting used by the other rows of the table (e.g. row 7):
#  - C: vertical-top alignment
#  - E: left/top alignment with wrapped text
$ws.Range("C9").VerticalAlignment = -4160
$ws.Range("E9").WrapText = $true
$ws.Range("E9").VerticalAlignment = -4160
$ws.Range("E9").HorizontalAlignment = -4131

# Row height for the newly filled (now 2-line) row
$ws.Rows("9").RowHeight = 30

# Widen column E to fit the longer text
$ws.Columns("E").ColumnWidth = 47.5

# Move/leave selection on A9 as in the final workbook
$ws.Range("A9").Select()
